$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Quantité Ingrédients" column (old column C)
# to hold a numeric "Quantite" value, shifting the old C:G columns to D:H.
$ws.Columns("C:C").Insert()

# Match the new column's width to column B's width (both are now 20.77734375 in the target file)
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Header
$ws.Range("C1").Value = "Quantite"

# Numeric quantity values (the leading number of each "Quantité Ingrédients" text)
$ws.Range("C2").Value = 300
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 30
$ws.Range("C9").Value = 30
$ws.Range("C10").Value = 50
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("C13").Value = 400
$ws.Range("C14").Value = 1

# Restore the selection/cursor state recorded in the saved file
$ws.Range("C15").Select() | Out-Null
